# edit.ps1 - apply "edit projects, temperature defn, KM's role" changes
$wb = $excel.ActiveWorkbook

# 1) Fix the "temperature" attribute definition on the ColumnHeadersNcp sheet
#    (row 7): it was mistakenly copy-pasted from the "salinity" row. Correct
#    the word "salinity" -> "temperature" in the definition text; unit stays
#    "celsius".
$ncp = $wb.Worksheets.Item("ColumnHeadersNcp")
$ncp.Range("B7").Value = "Underway thermosalinograph temperature in degrees Celsius. URI http://vocab.nerc.ac.uk/collection/P01/current/TEMPSZ01/"

# 2) Add Kate Morkeski as a new Personnel row (metadata Provider) on the
#    Personnel sheet.
$personnel = $wb.Worksheets.Item("Personnel")
$personnel.Range("A11").Value = "Kate"
$personnel.Range("C11").Value = "Morkeski"
$personnel.Range("D11").Value = "Northeast U.S. Shelf LTER"
$personnel.Range("E11").Value = "kmorkeski@whoi.edu"
$personnel.Range("F11").Value = "0000-0002-2903-5851"
$personnel.Range("G11").Value = "metadata Provider"
$personnel.Range("H11").Value = "Northeast U.S. Shelf LTER"
$personnel.Range("I11").Value = "NSF"
$personnel.Range("J11").Value = "OCE-2322676"

# Select the newly added row, matching the author's final selection state.
$personnel.Range("A11:J11").Select() | Out-Null
